$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.794.80"
$ws.Range("E2").Value = "  +1.17%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.540.65"
$ws.Range("E3").Value = "  +0.12%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "611.99"
$ws.Range("E5").Value = "  +0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "152.29"
$ws.Range("E6").Value = "  -1.38%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.539.86"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.482"
$ws.Range("E9").Value = "  -0.67%  "

$ws.Range("E10").Value = "  -1.35%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.08"
$ws.Range("E11").Value = "  +3.42%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.426"
$ws.Range("E12").Value = "  -1.00%  "

$ws.Range("E13").Value = "  -1.31%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.141.13"
$ws.Range("E14").Value = "  +0.17%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "31.95"
$ws.Range("E15").Value = "  -0.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.558.28"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.559.44"
$ws.Range("E17").Value = "  +0.86%  "

$ws.Range("E18").Value = "  -0.77%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("E19").Value = "  +0.09%  "

$ws.Range("E20").Value = "  -1.66%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.70"
$ws.Range("E21").Value = "  +3.41%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.30"
$ws.Range("E22").Value = "  -1.35%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.622"
$ws.Range("E23").Value = "  -2.63%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "76.96"
$ws.Range("E24").Value = "  -2.48%  "

$ws.Range("E25").Value = "  +4.79%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.683.27"
$ws.Range("E26").Value = "  +0.22%  "

$ws.Range("E27").Value = "  +0.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.18"
$ws.Range("E28").Value = "  -0.62%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.63"
$ws.Range("E29").Value = "  +3.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.54"
$ws.Range("E30").Value = "  -0.37%  "

$ws.Range("E31").Value = "  -4.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.169"
$ws.Range("E32").Value = "  +6.74%  "

$ws.Range("E33").Value = "  +0.19%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.67"
$ws.Range("E34").Value = "  -1.07%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.20"
$ws.Range("E35").Value = "  +0.10%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.528.19"
$ws.Range("E36").Value = "  -0.10%  "

$ws.Range("E37").Value = "  -2.78%  "

$ws.Range("E38").Value = "  -1.07%  "

$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "176.23"
$ws.Range("E41").Value = "  +0.14%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.20"
$ws.Range("E42").Value = "  +2.46%  "

$ws.Range("E43").Value = "  +2.10%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.39"
$ws.Range("E44").Value = "  -3.97%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.888"
$ws.Range("E45").Value = "  -0.53%  "

$ws.Range("E46").Value = "  +1.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.54"
$ws.Range("E47").Value = "  -0.68%  "

$ws.Range("E48").Value = "  -0.28%  "

$ws.Range("E49").Value = "  +3.39%  "

$ws.Range("E50").Value = "  -0.80%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.252"
$ws.Range("E51").Value = "  +0.22%  "
